# Updated cryptos list data (price + 1h volume change) pulled from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "39.944.74"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.203.15"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.31"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "86.93"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.511"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("E10").Value = "  +6.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.53"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.114"
$ws.Range("E13").Value = "  +2.63%  "
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.543.89"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.76"
$ws.Range("E16").Value = "  -1.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.202.05"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.732"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "39.861.08"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0884"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("E21").Value = "  -3.15%  "
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.42"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.93"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.00"
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.28"
$ws.Range("E31").Value = "  +2.28%  "
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.93"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("E35").Value = "  +5.23%  "
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.33"
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("E38").Value = "  +0.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0995"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.46"
$ws.Range("E41").Value = "  -2.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.072.10"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.14"
$ws.Range("E44").Value = "  +7.18%  "
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.85"
$ws.Range("E46").Value = "  -1.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.73"
$ws.Range("E47").Value = "  +2.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.93"
$ws.Range("E48").Value = "  -9.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.418.87"
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("E51").Value = "  +0.47%  "
